$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.05038766666666666
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.754144
$ws.Range("N2").Value = 11.262432
$ws.Range("O2").Value = 0.2855964853960988
$ws.Range("P2").Value = 0.2855964853960988
$ws.Range("Q2").Value = 0.1891625564906667
$ws.Range("R2").Value = 1.702463008416
$ws.Range("S2").Value = 0.2855964853960988
$ws.Range("T2").Value = 0.2855964853960988

# Row 3
$ws.Range("G3").Value = 0.05038766666666666
$ws.Range("O3").Value = 0.3020986986332443
$ws.Range("P3").Value = 0.3020986986332442
$ws.Range("S3").Value = 0.3020986986332443
$ws.Range("T3").Value = 0.3020986986332442

# Row 4
$ws.Range("G4").Value = 0.05038766666666666
$ws.Range("M4").Value = 5.376329333333334
$ws.Range("N4").Value = 16.128988
$ws.Range("O4").Value = 0.4090042262449046
$ws.Range("P4").Value = 0.4090042262449045
$ws.Range("Q4").Value = 0.2709006903382222
$ws.Range("R4").Value = 2.438106213044
$ws.Range("S4").Value = 0.4090042262449046
$ws.Range("T4").Value = 0.4090042262449045

# Row 5
$ws.Range("G5").Value = 0.05038766666666666
$ws.Range("M5").Value = 0.043386
$ws.Range("N5").Value = 0.130158
$ws.Range("O5").Value = 0.003300589725752433
$ws.Range("P5").Value = 0.003300589725752433
$ws.Range("Q5").Value = 0.002186119306
$ws.Range("R5").Value = 0.019675073754
$ws.Range("S5").Value = 0.003300589725752433
$ws.Range("T5").Value = 0.003300589725752433
